# Test <Listener and Screenshot capturing
#
# Adds a new "OpenAccount" worksheet (account-opening test data, with the
# customer name pulled back from AddCustomerTest via formula) and refreshes
# AddCustomerTest's sample rows 3-7 with real test data instead of the old
# chk/chk1/chk22 placeholders.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("AddCustomerTest")

# --- Refresh AddCustomerTest sample rows (3-7) -----------------------------
$ws1.Range("A3").Value = "Clinda "
$ws1.Range("B3").Value = "T"
$ws1.Range("C3").Value = 2222
$ws1.Range("D3").Value = "Customer added successfully"

$ws1.Range("A4").Value = "Bismi"
$ws1.Range("B4").Value = "S"
$ws1.Range("C4").Value = 222
$ws1.Range("D4").Value = "Customer added successfully"

$ws1.Range("A5").Value = "jira"
$ws1.Range("B5").Value = "R"
$ws1.Range("C5").Value = "1w2"
$ws1.Range("D5").Value = "Customer added successfully"

$ws1.Range("A6").Value = "Test"
$ws1.Range("B6").Value = "test2"
$ws1.Range("C6").Value = 11
$ws1.Range("D6").Value = "Customer added successfully"

$ws1.Range("A7").Value = "test"
$ws1.Range("B7").Value = "test3"
$ws1.Range("C7").Value = 33
$ws1.Range("D7").Value = "Customer added successfully"

# --- Add the OpenAccount worksheet right after AddCustomerTest ------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "OpenAccount"

$ws2.Columns.Item(1).ColumnWidth = 17.5703125
$ws2.Columns.Item(3).ColumnWidth = 27.28515625
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

$ws2.Range("A1").Value = "Customer"
$ws2.Range("B1").Value = "Currency"
$ws2.Range("C1").Value = "Message"

$ws2.Range("A2").Formula = "=AddCustomerTest!A2:B2"
$ws2.Range("B2").Value = "Pound"
$ws2.Range("C2").Value = "Account created successfully"

$ws2.Range("A3").Formula = "=AddCustomerTest!A3:B3"
$ws2.Range("B3").Value = "Dollar"
$ws2.Range("C3").Value = "Account created successfully"

$ws2.Range("A4").Formula = "=AddCustomerTest!A4:B4"
$ws2.Range("B4").Value = "Pound"
$ws2.Range("C4").Value = "Account created successfully"

$ws2.Range("A5").Formula = "=AddCustomerTest!A5:B5"
$ws2.Range("B5").Value = "Rupee"
$ws2.Range("C5").Value = "Account created successfully"

$ws2.Range("A6").Formula = "=AddCustomerTest!A6:B6"
$ws2.Range("B6").Value = "Rupee"
$ws2.Range("C6").Value = "Account created successfully"

$ws2.Range("A7").Formula = "=AddCustomerTest!A7:B7"
$ws2.Range("B7").Value = "Rupee"
$ws2.Range("C7").Value = "Account created successfully"

# --- Selections (OpenAccount ends up the active/visible tab) --------------
[void]$ws1.Range("A6").Select()
[void]$ws2.Range("A4:A5").Select()
